# API DOC small update
# - rename B13 path ("정해주세요" -> "recommendations/")
# - rename B18 path ("<int:articlePk>/" -> "<int:articlePk>/comments/")
# - add a new "back" / "v" / "todo" status column (G)
# - move active selection to E14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- text content updates -------------------------------------------------
$ws.Range("B13").Value = "recommendations/"
$ws.Range("B18").Value = "<int:articlePk>/comments/"

# --- new column G header ("back") ------------------------------------------
$ws.Range("C5").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = "back"
$ws.Range("G3").Interior.Color = 13434879
$ws.Range("G3").HorizontalAlignment = -4108
$ws.Range("G3").Borders.Item(7).Color = 11711154

# --- new column G body cells -------------------------------------------
$bodyRows = 4,5,6,7,8,9,10,11,12,14,15,16,17
foreach ($r in $bodyRows) {
    $ws.Range("C5").Copy()
    $target = $ws.Range("G" + $r)
    $target.PasteSpecial(-4122)
    $target.Value = "v"
    $target.HorizontalAlignment = -4108
}

# rows 13, 18, 19, 20 use the plain default font style
$ws.Range("G13").Value = "todo"
$ws.Range("G13").HorizontalAlignment = -4108

$ws.Range("G18").Value = "v"
$ws.Range("G18").HorizontalAlignment = -4108

$ws.Range("G19").Value = "v"
$ws.Range("G19").HorizontalAlignment = -4108

$ws.Range("G20").Value = "v"
$ws.Range("G20").HorizontalAlignment = -4108

# --- selection -------------------------------------------------------------
$ws.Range("E14").Select()
